$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36 by copying row 35 (keeps formulas/format consistent
# with the rest of the t-axis constraint table) and inserting the copy,
# which shifts the old rows 36/37 down to 37/38.
$ws.Rows.Item(35).Copy()
$ws.Rows.Item(36).Insert()

# Row 34 ("Both TOL inc") - prog2dams/prog2off constraint now also applies
$ws.Range("C34").Value = $true

# Row 35 becomes the new "TOL 1 only" constraint row
$ws.Range("C35").Value = $true
$ws.Range("D35").Value = "TOL 1 only"
$ws.Range("E35").Value = $true
$ws.Range("F35").Value = $false

# Row 36 is the newly inserted "TOL 2 only" constraint row
$ws.Range("B36").Value = $false
$ws.Range("C36").Value = $true
$ws.Range("D36").Value = "TOL 2 only"
$ws.Range("E36").Value = $false
$ws.Range("F36").Value = $true

# Row 37 ("BBT inc", previously row 36) also gets the new constraint
$ws.Range("C37").Value = $true
